# Update the "poisson" data sheet for "semana 40 de 2025"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (evento 113)
$ws.Range("C3").Value = 6
$ws.Range("E3").Value = 0.01

# Row 4 (evento 115)
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1

# Row 5 (evento 155)
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 0.09

# Row 6 (evento 210)
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 0.27

# Row 7 (evento 215)
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0.18

# Row 9 (evento 300)
$ws.Range("C9").Value = 38
$ws.Range("D9").Value = 47
$ws.Range("E9").Value = 0.02

# Row 11 (evento 340)
$ws.Range("C11").Value = 2
$ws.Range("E11").Value = 0.09

# Row 12 (evento 342)
$ws.Range("C12").Value = 7
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 0.09

# Row 13 (evento 346)
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 0

# Row 17 (evento 356)
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = 0.09

# Row 19 (evento 365)
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 0.08

# Row 22 (evento 455)
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 0.18

# Row 25 (evento 549)
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 0.15

# Row 26 (evento 560)
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0.27

# Row 29 (evento 620)
$ws.Range("C29").Value = 2
$ws.Range("E29").Value = 0.14

# Row 31 (evento 750)
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 0.27

# Row 33 (evento 813)
$ws.Range("C33").Value = 7
$ws.Range("D33").Value = 5
$ws.Range("E33").Value = 0.13

# Row 34 (evento 831)
$ws.Range("C34").Value = 10
$ws.Range("D34").Value = 4
$ws.Range("E34").Value = 0.02

# Row 35 (evento 850)
$ws.Range("C35").Value = 10
$ws.Range("D35").Value = 5
$ws.Range("E35").Value = 0.04
